$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-9 from 45204 (2023-10-05)
# to 45207 (2023-10-08), matching the automatic update reflected in the diff.
$ws.Range("C2:C9").Value = 45207
